$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, pushing existing rows 18-25 down to 19-26
$ws.Rows.Item(18).Insert()

# Fill in the new row 18 with the new weekly record
$ws.Cells.Item(18, 1).Value = 7
$ws.Cells.Item(18, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(18, 3).Value = "Ñuble"
$ws.Cells.Item(18, 4).Value = 44966
$ws.Cells.Item(18, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18, 5).Value = 16
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100101
$ws.Cells.Item(18, 8).Value = "Berries"
$ws.Cells.Item(18, 9).Value = 100101001
$ws.Cells.Item(18, 10).Value = "Arándano (blue)"
$ws.Cells.Item(18, 11).Value = "Sin especificar"
$ws.Cells.Item(18, 12).Value = "Segunda"
$ws.Cells.Item(18, 13).Value = 30
$ws.Cells.Item(18, 14).Value = 2500
$ws.Cells.Item(18, 15).Value = 2500
$ws.Cells.Item(18, 16).Value = 2500
$ws.Cells.Item(18, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(18, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(18, 19).Value = 1250
$ws.Cells.Item(18, 20).Value = 2
